$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the covariance matrix in G2:J5 using COVARIANCE.S, referencing
# the data columns A:D (each row of the matrix keeps the first argument's
# range absolute while the second argument floats across A2:A11..D2:D11).

$ws.Range("G2").Formula = "=_xlfn.COVARIANCE.S(`$A`$2:`$A`$11,A2:A11)"
$ws.Range("H2").Formula = "=_xlfn.COVARIANCE.S(`$A`$2:`$A`$11,B2:B11)"
$ws.Range("I2").Formula = "=_xlfn.COVARIANCE.S(`$A`$2:`$A`$11,C2:C11)"
$ws.Range("J2").Formula = "=_xlfn.COVARIANCE.S(`$A`$2:`$A`$11,D2:D11)"

$ws.Range("G3").Formula = "=_xlfn.COVARIANCE.S(`$B`$2:`$B`$11,A2:A11)"
$ws.Range("H3").Formula = "=_xlfn.COVARIANCE.S(`$B`$2:`$B`$11,B2:B11)"
$ws.Range("I3").Formula = "=_xlfn.COVARIANCE.S(`$B`$2:`$B`$11,C2:C11)"
$ws.Range("J3").Formula = "=_xlfn.COVARIANCE.S(`$B`$2:`$B`$11,D2:D11)"

$ws.Range("G4").Formula = "=_xlfn.COVARIANCE.S(`$C`$2:`$C`$11,A2:A11)"
$ws.Range("H4").Formula = "=_xlfn.COVARIANCE.S(`$C`$2:`$C`$11,B2:B11)"
$ws.Range("I4").Formula = "=_xlfn.COVARIANCE.S(`$C`$2:`$C`$11,C2:C11)"
$ws.Range("J4").Formula = "=_xlfn.COVARIANCE.S(`$C`$2:`$C`$11,D2:D11)"

$ws.Range("G5").Formula = "=_xlfn.COVARIANCE.S(`$D`$2:`$D`$11,A2:A11)"
$ws.Range("H5").Formula = "=_xlfn.COVARIANCE.S(`$D`$2:`$D`$11,B2:B11)"
$ws.Range("I5").Formula = "=_xlfn.COVARIANCE.S(`$D`$2:`$D`$11,C2:C11)"
$ws.Range("J5").Formula = "=_xlfn.COVARIANCE.S(`$D`$2:`$D`$11,D2:D11)"

# Update sheet view to match final state (zoom + selected cell)
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("I7").Select()
